$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new student block (rows 26-29) for "S Dey", roll no. 7,
# mirroring the 4-row-per-student layout already used in the sheet.

$ws.Cells.Item(26, 1).Value = 7
$ws.Cells.Item(26, 2).Value = "S Dey"
$ws.Cells.Item(26, 3).Value = "First Summative"
$ws.Cells.Item(26, 4).Value = 9
$ws.Cells.Item(26, 5).Value = 11
$ws.Cells.Item(26, 6).Value = 13
$ws.Cells.Item(26, 7).Value = 18
$ws.Cells.Item(26, 8).Value = 36
$ws.Cells.Item(26, 9).Value = 30
$ws.Cells.Item(26, 10).Value = 29
$ws.Cells.Item(26, 11).Value = 27

$ws.Cells.Item(27, 3).Value = "Second Summative"
$ws.Cells.Item(27, 4).Value = 10
$ws.Cells.Item(27, 5).Value = 11
$ws.Cells.Item(27, 6).Value = 8
$ws.Cells.Item(27, 7).Value = 31
$ws.Cells.Item(27, 8).Value = 50
$ws.Cells.Item(27, 9).Value = 31
$ws.Cells.Item(27, 10).Value = 27
$ws.Cells.Item(27, 11).Value = 31

$ws.Cells.Item(28, 3).Value = "Third Summative Theory"
$ws.Cells.Item(28, 4).Value = 14
$ws.Cells.Item(28, 5).Value = 18
$ws.Cells.Item(28, 6).Value = 12
$ws.Cells.Item(28, 7).Value = 51
$ws.Cells.Item(28, 8).Value = 90
$ws.Cells.Item(28, 9).Value = 51
$ws.Cells.Item(28, 10).Value = 51
$ws.Cells.Item(28, 11).Value = 51

$ws.Cells.Item(29, 3).Value = "Third Summative Practical"
$ws.Cells.Item(29, 4).Value = 9
$ws.Cells.Item(29, 5).Value = 8
$ws.Cells.Item(29, 6).Value = 7
$ws.Cells.Item(29, 7).Value = 9
$ws.Cells.Item(29, 8).Value = 10
$ws.Cells.Item(29, 9).Value = 9
$ws.Cells.Item(29, 10).Value = 9
$ws.Cells.Item(29, 11).Value = 9

# Update selection to match the authored state (G30, empty cell just
# below the newly added block).
$ws.Range("G30").Select()
